# Delete row 591 ("「ラヤンはねむれない」...") entirely, shifting all
# subsequent rows up by one. This matches the commit's removal of that
# post entry from the posts sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(591).Delete()
